$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: run a Find/Replace over the whole document that merges adjacent
# runs containing the same plain text into a single run. Word (and this
# interop layer) collapses a matched range that spans several runs into one
# run when its .Text is replaced, which is exactly the "run merge" behaviour
# shown throughout the diff.
# ---------------------------------------------------------------------------
function Merge-Text([string]$text) {
    $d.Content.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, $text, 2) | Out-Null
}

# 2º paragraph: " Após a instalação criar uma pasta no "+"seu computador"+" que será utilizada como repositório"+" local"
Merge-Text(" Após a instalação criar uma pasta no seu computador que será utilizada como repositório local")

# 3º paragraph: " "+"Cl"+"ica com o botão direito do mouse e escolher a opção: "
Merge-Text(" Clica com o botão direito do mouse e escolher a opção: ")

# 4º paragraph: " Com o prompt de comando do git aberto, digite o comando:"+" "
Merge-Text(" Com o prompt de comando do git aberto, digite o comando: ")

# 5º paragraph: " Adicionar a essa pasta o arquivo que será utilizado para criar um"+"a"+" versão junto ao seu repositório"
Merge-Text(" Adicionar a essa pasta o arquivo que será utilizado para criar uma versão junto ao seu repositório")

# 6º paragraph: " Você verá o "+"arquivo"+" sinalizado em "
Merge-Text(" Você verá o arquivo sinalizado em ")

# 6º paragraph (continued): ", pois "+"ele"+" ainda não está pronto para ser "
Merge-Text(", pois ele ainda não está pronto para ser ")

# 7º paragraph: " "+" "+". Com isso, não só o arquivo que está sinalizado em "
Merge-Text("  . Com isso, não só o arquivo que está sinalizado em ")

# ---------------------------------------------------------------------------
# 12º paragraph: "git push --set-upstream origin master" becomes two runs:
# "git push --set-upstream origin " and "receiving" (same bold/size 28
# formatting), and the later standalone "branch master" reference becomes
# "branch receiving".
# ---------------------------------------------------------------------------
$phrase = $d.Content
$phrase.Find.Execute("git push --set-upstream origin master") | Out-Null

# Locate just the "master" word inside the matched phrase so only it gets
# replaced (leaving "git push --set-upstream origin " intact).
$rMaster = $d.Range($phrase.Start, $phrase.End)
$rMaster.Find.Execute("master") | Out-Null
$masterStart = $rMaster.Start
$rMaster.Text = "receiving"
$newEnd = $masterStart + 9
$r2 = $d.Range($masterStart, $newEnd)
# Re-assert (nudge) the run formatting so the interop layer keeps this
# fragment as its own run instead of re-merging it with its bold/size-28
# neighbour that now has identical formatting.
$r2.Font.Size = 14.5
$r2.Font.Size = 14
$r2.Font.Bold = $true

# Replace the remaining standalone "master" (the "branch master" mention),
# scoped to the current paragraph only so no other "master" occurrence
# later in the document is touched.
$para = $r2.Paragraphs(1)
$paraRange = $d.Range($para.Range.Start, $para.Range.End)
$paraRange.Find.Execute("master", $true, $false, $false, $false, $false, $true, 1, $false, "receiving", 2) | Out-Null
